$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest crypto snapshot data.
# Some Price values (e.g. "614.04") look like plain numbers, so Excel would silently
# convert them to numeric cells (losing formatting such as trailing zeros) unless the
# cell is explicitly forced to Text format first. We do this only where needed and
# restore the original (default) cell style afterwards so no visible formatting changes.

$ws.Range("D2").Value = '70.922.68'
$ws.Range("E2").Value = '  +2.77%  '
$ws.Range("D3").Value = '3.566.97'
$ws.Range("E3").Value = '  +2.09%  '
$ws.Range("E4").Value = '  -0.01%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '614.04'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +6.44%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '172.22'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.69%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.619'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +2.58%  '
$ws.Range("D8").Value = '3.561.90'
$ws.Range("E8").Value = '  +2.14%  '
$ws.Range("E9").Value = '  -0.02%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.198'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +5.54%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '7.37'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +13.87%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.589'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +1.70%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '46.78'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("D15").Value = '4.135.02'
$ws.Range("E15").Value = '  +1.92%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '8.42'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -0.83%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '619.45'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("D18").Value = '70.924.24'
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("D19").Value = '3.557.66'
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("E20").Value = '  -1.55%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '17.41'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +0.96%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.884'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.46%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '9.49'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -14.26%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '15.80'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '97.01'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  +1.75%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +0.18%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '33.65'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +3.59%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '9.13'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -1.55%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '8.56'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +1.31%  '
$ws.Range("E32").Value = '  -2.07%  '
$ws.Range("E33").Value = '  -0.16%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '6.97'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.69%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '572.37'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -8.70%  '
$ws.Range("E36").Value = '  +6.24%  '
$ws.Range("E37").Value = '  -0.64%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '10.88'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +1.82%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '57.67'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +2.15%  '
$ws.Range("E40").Value = '  +7.14%  '
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("E42").Value = '  +5.83%  '
$ws.Range("D43").Value = '3.376.44'
$ws.Range("E43").Value = '  +0.82%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.321'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -1.45%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '3.00'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +8.61%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '33.05'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("D47").Value = '0.0₃0706'
$ws.Range("E47").Value = '  +2.36%  '
$ws.Range("E48").Value = '  +2.98%  '
$ws.Range("E49").Value = '  +1.23%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '133.71'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +1.25%  '
$ws.Range("E51").Value = '  +1.78%  '
